$wb = $excel.ActiveWorkbook

# --- Refresh standings (Total/Sofrido/Saldo Cartola) for every group sheet ---
# Values below are the latest scrape output from the Libertadores dataset notebook run.
$standings = @{
    "Grupo A" = @{
        2 = @{ G = 66.86; H = 46.79; I = 20.07 }
        3 = @{ G = 48.5; H = 43.56; I = 4.939999999999998 }
        4 = @{ G = 46.79; H = 66.86; I = -20.07 }
        5 = @{ G = 43.56; H = 48.5; I = -4.939999999999998 }
    }
    "Grupo B" = @{
        2 = @{ G = 65.06; H = 47.16; I = 17.90000000000001 }
        3 = @{ G = 59.69; H = 43.56; I = 16.13 }
        4 = @{ G = 47.16; H = 65.06; I = -17.90000000000001 }
        5 = @{ G = 43.56; H = 59.69; I = -16.13 }
    }
    "Grupo C" = @{
        2 = @{ G = 71.46; H = 57.26; I = 14.2 }
        3 = @{ G = 71.36; H = 34.36; I = 37 }
        4 = @{ G = 57.26; H = 71.46; I = -14.2 }
        5 = @{ G = 34.36; H = 71.36; I = -37 }
    }
    "Grupo D" = @{
        2 = @{ G = 71.36; H = 49.36; I = 22 }
        3 = @{ G = 63.76; H = 50.69; I = 13.07 }
        4 = @{ G = 50.69; H = 63.76; I = -13.07 }
        5 = @{ G = 49.36; H = 71.36; I = -22 }
    }
    "Grupo E" = @{
        2 = @{ G = 81.76; H = 56.09; I = 25.67 }
        3 = @{ G = 59.65; H = 50.76; I = 8.89 }
        4 = @{ G = 56.09; H = 81.76; I = -25.67 }
        5 = @{ G = 50.76; H = 59.65; I = -8.89 }
    }
    "Grupo F" = @{
        2 = @{ G = 58.26; H = 33.96; I = 24.3 }
        3 = @{ G = 54.95; H = 52.39; I = 2.560000000000002 }
        4 = @{ G = 52.39; H = 54.95; I = -2.560000000000002 }
        5 = @{ G = 33.96; H = 58.26; I = -24.3 }
    }
    "Grupo G" = @{
        2 = @{ G = 71.16; H = 48.29; I = 22.87 }
        3 = @{ G = 59.36; H = 48.89; I = 10.47 }
        4 = @{ G = 48.89; H = 59.36; I = -10.47 }
        5 = @{ G = 48.29; H = 71.16; I = -22.87 }
    }
    "Grupo H" = @{
        2 = @{ G = 63.76; H = 54.66; I = 9.100000000000001 }
        3 = @{ G = 56.65; H = 38.66; I = 17.99 }
        4 = @{ G = 54.66; H = 63.76; I = -9.100000000000001 }
        5 = @{ G = 38.66; H = 56.65; I = -17.99 }
    }
}

foreach ($sheetName in $standings.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $standings[$sheetName]
    foreach ($r in $rows.Keys) {
        $vals = $rows[$r]
        $ws.Range("G$r").Value2 = $vals.G
        $ws.Range("H$r").Value2 = $vals.H
        $ws.Range("I$r").Value2 = $vals.I
    }
}

# --- Fix team-name ordering that drifted in the raw scrape (stale shared-string refs) ---
$wsA = $wb.Worksheets.Item("Grupo A")
$tmp = $wsA.Range("B3").Value2
$wsA.Range("B3").Value2 = $wsA.Range("B5").Value2
$wsA.Range("B5").Value2 = $tmp

$wsC = $wb.Worksheets.Item("Grupo C")
$tmp = $wsC.Range("B2").Value2
$wsC.Range("B2").Value2 = $wsC.Range("B3").Value2
$wsC.Range("B3").Value2 = $tmp

$wsD = $wb.Worksheets.Item("Grupo D")
$tmp = $wsD.Range("B4").Value2
$wsD.Range("B4").Value2 = $wsD.Range("B5").Value2
$wsD.Range("B5").Value2 = $tmp

